{"js": "// Replace the title/intro, the \"What we like\" / \"What we don't like\" bullet\n// text, and the closing bold/italic summary lines per the commit diff.\nconst replacements = [\n  [\n    \"Play Free Disco Diamonds Slot - A Retro Disco Experience\",\n    \"Play Disco Diamonds Free - Captivating Gameplay with Mesmerizing Visuals\",\n  ],\n  [\n    \"Two special modes offer immediate payouts, free spins, and multipliers.\",\n    \"Captivating gameplay with special modes\",\n  ],\n  [\n    \"Fascinating neon-colored visuals combined with retro-style graphics.\",\n    \"Mesmerizing visuals with neon colors\",\n  ],\n  [\n    \"Classic and entertaining gameplay.\",\n    \"Unique charm and disco aesthetics\",\n  ],\n  [\n    \"Disco Diamonds offers its players a high level of customizability.\",\n    \"Variable RTP for an overall satisfying experience\",\n  ],\n  [\n    \"Disco Diamonds has a higher volatility than some other slot games.\",\n    \"High volatility may not suit all players\",\n  ],\n  [\n    \"A lack of progressive jackpot might deter enthusiasts.\",\n    \"Limited selection of similar games\",\n  ],\n  [\n    \"Read our review of the Disco Diamonds slot game. Dance to the rhythm with two special modes that offer payouts, free spins, and multipliers. Play for free now.\",\n    \"Read our review of Disco Diamonds and play this captivating game for free. Mesmerizing visuals and unique charm await!\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the title/intro, the \"What we like\" / \"What we don't like\" bullet\n# text, and the closing bold/italic summary lines per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Free Disco Diamonds Slot - A Retro Disco Experience\", \"Play Disco Diamonds Free - Captivating Gameplay with Mesmerizing Visuals\"),\n    @(\"Two special modes offer immediate payouts, free spins, and multipliers.\", \"Captivating gameplay with special modes\"),\n    @(\"Fascinating neon-colored visuals combined with retro-style graphics.\", \"Mesmerizing visuals with neon colors\"),\n    @(\"Classic and entertaining gameplay.\", \"Unique charm and disco aesthetics\"),\n    @(\"Disco Diamonds offers its players a high level of customizability.\", \"Variable RTP for an overall satisfying experience\"),\n    @(\"Disco Diamonds has a higher volatility than some other slot games.\", \"High volatility may not suit all players\"),\n    @(\"A lack of progressive jackpot might deter enthusiasts.\", \"Limited selection of similar games\"),\n    @(\"Read our review of the Disco Diamonds slot game. Dance to the rhythm with two special modes that offer payouts, free spins, and multipliers. Play for free now.\", \"Read our review of Disco Diamonds and play this captivating game for free. Mesmerizing visuals and unique charm await!\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
